$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'281.36"
$ws.Range("D4").Value = "'6.237"
$ws.Range("D5").Value = "'0.06140"
$ws.Range("D7").Value = "'6.560"
$ws.Range("D8").Value = "'1.455"
$ws.Range("D11").Value = "'0.1631"
$ws.Range("D12").Value = "'0.08304"
$ws.Range("D13").Value = "'0.03536"
$ws.Range("D14").Value = "'0.03210"
$ws.Range("D15").Value = "'0.09142"
$ws.Range("D16").Value = "'3.731"
$ws.Range("D17").Value = "'0.001645"
$ws.Range("D18").Value = "'0.04635"
$ws.Range("D19").Value = "'0.006468"
$ws.Range("D20").Value = "'0.006163"
$ws.Range("D23").Value = "'3.804"
$ws.Range("D25").Value = "'0.3338"
$ws.Range("D40").Value = "'0.04657"
$ws.Range("D41").Value = "'0.006302"
$ws.Range("D42").Value = "'0.007176"
$ws.Range("D43").Value = "'0.1097"
$ws.Range("D44").Value = "'0.01134"
$ws.Range("D45").Value = "'0.00006427"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.002942"
$ws.Range("D49").Value = "'0.00001901"
$ws.Range("D50").Value = "'0.01240"
